# Ticket 171: ajustado el nombre del archivo de incremento de claro, segun comentario 579
#
# Cell A4 on "Hoja1" holds a rich-text placeholder made of two runs:
#   run 1: a leading BOM (zero-width) character in ".AppleSystemUIFont" (sz 3.9)
#   run 2: the placeholder text "${ID_ORDEN_SERVICIO}" in Arial (sz 10)
#
# Per comment 579 on ticket 171, run 2's text must become "${ID_SERVICIO} "
# (note the trailing space), keeping each run's original font.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws3 = $wb.Worksheets.Item("Hoja3")

$cell = $ws1.Range("A4")

# How many characters currently make up the cell (BOM + placeholder text).
$totalLen = $cell.Characters().Count

# Replace everything after the leading BOM character (position 1) with the
# new placeholder text.
$target = $cell.Characters(2, $totalLen - 1)
$target.Text = '${ID_SERVICIO} '

# Re-apply the original formatting to each run, since replacing the text
# resets character-level font info on the cell.
$newTotalLen = $cell.Characters().Count

$bom = $cell.Characters(1, 1)
$bom.Font.Name = ".AppleSystemUIFont"
$bom.Font.Size = 3.9

$rest = $cell.Characters(2, $newTotalLen - 1)
$rest.Font.Name = "Arial"
$rest.Font.Size = 10

# Hoja2 and Hoja3 had a stray multi-range selection (A4 and A1) left over
# from editing; reset each sheet's selection back to just A1, then restore
# the original active sheet/selection (A4 on Hoja1).
$ws2.Select() | Out-Null
$ws2.Range("A1").Select() | Out-Null

$ws3.Select() | Out-Null
$ws3.Range("A1").Select() | Out-Null

$ws1.Select() | Out-Null
$ws1.Range("A4").Select() | Out-Null
